$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# PHASE 1: formatting (fill/font). Order matters for style-index allocation:
# the bold+yellow combo (header/level cells, column A) must be created before
# the plain yellow-only combo (the rest of the level row), so that they land
# on cellXfs index 2 and 3 respectively.
# ---------------------------------------------------------------------------
$levelHeaderRows = 8, 11, 14
foreach ($r in $levelHeaderRows) {
    $ws.Range("A$r").Font.Bold = $true
    $ws.Range("A$r").Interior.Color = 65535
}
foreach ($r in $levelHeaderRows) {
    $ws.Range("B$r").Interior.Color = 65535
    $ws.Range("C$r").Interior.Color = 65535
    $ws.Range("D$r").Interior.Color = 65535
    $ws.Range("E$r").Interior.Color = 65535
}

# ---------------------------------------------------------------------------
# PHASE 2: cell values, in the exact order the new shared strings were first
# introduced so that the rebuilt sharedStrings table lines up.
# ---------------------------------------------------------------------------

# --- Level 1: Address -------------------------------------------------
$ws.Range("C8").Value = "LVL"
$ws.Range("D8").Value = 1

$ws.Range("A8").Value = "Address"

$ws.Range("A9").Value = "AddressId"
$ws.Range("C9").Value = "PK"
$ws.Range("D9").Value = "Num"
$ws.Range("E9").Value = 4

$ws.Range("A10").Value = "AddressName"
$ws.Range("D10").Value = "Char"
$ws.Range("E10").Value = 20

# --- Level 2: Properties ------------------------------------------------
$ws.Range("A11").Value = "Properties"
$ws.Range("C11").Value = "LVL"
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 1

$ws.Range("A12").Value = "PropertyId"
$ws.Range("C12").Value = "PK"
$ws.Range("D12").Value = "Num"
$ws.Range("E12").Value = 4

$ws.Range("A13").Value = "PropertyName"
$ws.Range("D13").Value = "Char"
$ws.Range("E13").Value = 20

# --- New AutoNumber column on the Customer table ------------------------
$ws.Range("G4").Font.Bold = $true
$ws.Range("G4").Value = "AutoNumber"
$ws.Range("G5").Value = $true

# --- Level 3: Features ----------------------------------------------------
$ws.Range("A14").Value = "Features"
$ws.Range("C14").Value = "LVL"
$ws.Range("D14").Value = 3

$ws.Range("A15").Value = "FeatureId"
$ws.Range("C15").Value = "PK"
$ws.Range("D15").Value = "Num"
$ws.Range("E15").Value = 4

$ws.Range("A16").Value = "FeatureName"
$ws.Range("D16").Value = "Char"
$ws.Range("E16").Value = 20

# ---------------------------------------------------------------------------
# PHASE 3: column width + selection, matching the sheetView / cols changes.
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 15.736979166666666
$ws.Range("G6").Select() | Out-Null

Write-Host "done"
